$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2025-12-03 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-12-04 Thursday", 2) | Out-Null

# Update each answer cell in the table by row/column position
$tbl = $d.Tables.Item(1)

$tbl.Rows.Item(1).Cells.Item(1).Range.Text = "13+56=69"
$tbl.Rows.Item(1).Cells.Item(2).Range.Text = "16+2=18"
$tbl.Rows.Item(1).Cells.Item(3).Range.Text = "98-58=40"
$tbl.Rows.Item(1).Cells.Item(4).Range.Text = "46-41=5"
$tbl.Rows.Item(1).Cells.Item(5).Range.Text = "12+32=44"
$tbl.Rows.Item(2).Cells.Item(1).Range.Text = "85-48=37"
$tbl.Rows.Item(2).Cells.Item(2).Range.Text = "91-47=44"
$tbl.Rows.Item(2).Cells.Item(3).Range.Text = "37-15=22"
$tbl.Rows.Item(2).Cells.Item(4).Range.Text = "27-0=27"
$tbl.Rows.Item(2).Cells.Item(5).Range.Text = "54-12=42"
$tbl.Rows.Item(3).Cells.Item(1).Range.Text = "47+37=84"
$tbl.Rows.Item(3).Cells.Item(2).Range.Text = "49-45=4"
$tbl.Rows.Item(3).Cells.Item(3).Range.Text = "36+8=44"
$tbl.Rows.Item(3).Cells.Item(4).Range.Text = "59-42=17"
$tbl.Rows.Item(3).Cells.Item(5).Range.Text = "76-48=28"
$tbl.Rows.Item(4).Cells.Item(1).Range.Text = "41+36=77"
$tbl.Rows.Item(4).Cells.Item(2).Range.Text = "66-17=49"
$tbl.Rows.Item(4).Cells.Item(3).Range.Text = "72-46=26"
$tbl.Rows.Item(4).Cells.Item(4).Range.Text = "17+21=38"
$tbl.Rows.Item(4).Cells.Item(5).Range.Text = "37+11=48"
$tbl.Rows.Item(5).Cells.Item(1).Range.Text = "98-24=74"
$tbl.Rows.Item(5).Cells.Item(2).Range.Text = "95-33=62"
$tbl.Rows.Item(5).Cells.Item(3).Range.Text = "74+6=80"
$tbl.Rows.Item(5).Cells.Item(4).Range.Text = "13+75=88"
$tbl.Rows.Item(5).Cells.Item(5).Range.Text = "0+76=76"
$tbl.Rows.Item(6).Cells.Item(1).Range.Text = "16+65=81"
$tbl.Rows.Item(6).Cells.Item(2).Range.Text = "50+1=51"
$tbl.Rows.Item(6).Cells.Item(3).Range.Text = "13+32=45"
$tbl.Rows.Item(6).Cells.Item(4).Range.Text = "31+56=87"
$tbl.Rows.Item(6).Cells.Item(5).Range.Text = "59-4=55"
$tbl.Rows.Item(7).Cells.Item(1).Range.Text = "34+12=46"
$tbl.Rows.Item(7).Cells.Item(2).Range.Text = "68+20=88"
$tbl.Rows.Item(7).Cells.Item(3).Range.Text = "62-62=0"
$tbl.Rows.Item(7).Cells.Item(4).Range.Text = "11+74=85"
$tbl.Rows.Item(7).Cells.Item(5).Range.Text = "92-44=48"
$tbl.Rows.Item(8).Cells.Item(1).Range.Text = "7+43=50"
$tbl.Rows.Item(8).Cells.Item(2).Range.Text = "72+23=95"
$tbl.Rows.Item(8).Cells.Item(3).Range.Text = "94-40=54"
$tbl.Rows.Item(8).Cells.Item(4).Range.Text = "54-42=12"
$tbl.Rows.Item(8).Cells.Item(5).Range.Text = "66+31=97"
$tbl.Rows.Item(9).Cells.Item(1).Range.Text = "24-12=12"
$tbl.Rows.Item(9).Cells.Item(2).Range.Text = "88-9=79"
$tbl.Rows.Item(9).Cells.Item(3).Range.Text = "74-45=29"
$tbl.Rows.Item(9).Cells.Item(4).Range.Text = "49-4=45"
$tbl.Rows.Item(9).Cells.Item(5).Range.Text = "49-3=46"
$tbl.Rows.Item(10).Cells.Item(1).Range.Text = "95-31=64"
$tbl.Rows.Item(10).Cells.Item(2).Range.Text = "7+15=22"
$tbl.Rows.Item(10).Cells.Item(3).Range.Text = "52-7=45"
$tbl.Rows.Item(10).Cells.Item(4).Range.Text = "88-30=58"
$tbl.Rows.Item(10).Cells.Item(5).Range.Text = "73-55=18"
$tbl.Rows.Item(11).Cells.Item(1).Range.Text = "77-37=40"
$tbl.Rows.Item(11).Cells.Item(2).Range.Text = "20+60=80"
$tbl.Rows.Item(11).Cells.Item(3).Range.Text = "84-39=45"
$tbl.Rows.Item(11).Cells.Item(4).Range.Text = "14+54=68"
$tbl.Rows.Item(11).Cells.Item(5).Range.Text = "36+26=62"
$tbl.Rows.Item(12).Cells.Item(1).Range.Text = "85-16=69"
$tbl.Rows.Item(12).Cells.Item(2).Range.Text = "0+87=87"
$tbl.Rows.Item(12).Cells.Item(3).Range.Text = "55+19=74"
$tbl.Rows.Item(12).Cells.Item(4).Range.Text = "96-10=86"
$tbl.Rows.Item(12).Cells.Item(5).Range.Text = "83+2=85"
$tbl.Rows.Item(13).Cells.Item(1).Range.Text = "90-47=43"
$tbl.Rows.Item(13).Cells.Item(2).Range.Text = "98-95=3"
$tbl.Rows.Item(13).Cells.Item(3).Range.Text = "8+19=27"
$tbl.Rows.Item(13).Cells.Item(4).Range.Text = "13+81=94"
$tbl.Rows.Item(13).Cells.Item(5).Range.Text = "29+55=84"
$tbl.Rows.Item(14).Cells.Item(1).Range.Text = "91-82=9"
$tbl.Rows.Item(14).Cells.Item(2).Range.Text = "29-20=9"
$tbl.Rows.Item(14).Cells.Item(3).Range.Text = "66-61=5"
$tbl.Rows.Item(14).Cells.Item(4).Range.Text = "67+30=97"
$tbl.Rows.Item(14).Cells.Item(5).Range.Text = "89-26=63"
$tbl.Rows.Item(15).Cells.Item(1).Range.Text = "82-70=12"
$tbl.Rows.Item(15).Cells.Item(2).Range.Text = "69-15=54"
$tbl.Rows.Item(15).Cells.Item(3).Range.Text = "52+30=82"
$tbl.Rows.Item(15).Cells.Item(4).Range.Text = "85+1=86"
$tbl.Rows.Item(15).Cells.Item(5).Range.Text = "48+8=56"
$tbl.Rows.Item(16).Cells.Item(1).Range.Text = "7+89=96"
$tbl.Rows.Item(16).Cells.Item(2).Range.Text = "94-10=84"
$tbl.Rows.Item(16).Cells.Item(3).Range.Text = "84-39=45"
$tbl.Rows.Item(16).Cells.Item(4).Range.Text = "17+74=91"
$tbl.Rows.Item(16).Cells.Item(5).Range.Text = "4+76=80"
$tbl.Rows.Item(17).Cells.Item(1).Range.Text = "93-74=19"
$tbl.Rows.Item(17).Cells.Item(2).Range.Text = "12+13=25"
$tbl.Rows.Item(17).Cells.Item(3).Range.Text = "90-0=90"
$tbl.Rows.Item(17).Cells.Item(4).Range.Text = "16+9=25"
$tbl.Rows.Item(17).Cells.Item(5).Range.Text = "9+62=71"
$tbl.Rows.Item(18).Cells.Item(1).Range.Text = "71+4=75"
$tbl.Rows.Item(18).Cells.Item(2).Range.Text = "93-38=55"
$tbl.Rows.Item(18).Cells.Item(3).Range.Text = "72+9=81"
$tbl.Rows.Item(18).Cells.Item(4).Range.Text = "33+4=37"
$tbl.Rows.Item(18).Cells.Item(5).Range.Text = "49+34=83"
$tbl.Rows.Item(19).Cells.Item(1).Range.Text = "68-62=6"
$tbl.Rows.Item(19).Cells.Item(2).Range.Text = "70-9=61"
$tbl.Rows.Item(19).Cells.Item(3).Range.Text = "37+8=45"
$tbl.Rows.Item(19).Cells.Item(4).Range.Text = "71-32=39"
$tbl.Rows.Item(19).Cells.Item(5).Range.Text = "77+1=78"
$tbl.Rows.Item(20).Cells.Item(1).Range.Text = "99-95=4"
$tbl.Rows.Item(20).Cells.Item(2).Range.Text = "7-7=0"
$tbl.Rows.Item(20).Cells.Item(3).Range.Text = "19+76=95"
$tbl.Rows.Item(20).Cells.Item(4).Range.Text = "39+56=95"
$tbl.Rows.Item(20).Cells.Item(5).Range.Text = "65+31=96"
